$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 2
$ws.Range("F4").Value = -2
$ws.Range("F7").Value = -4
$ws.Range("F9").Value = -9
$ws.Range("F13").Value = -4
$ws.Range("F18").Value = 0
$ws.Range("F20").Value = -3
$ws.Range("F23").Value = -2
$ws.Range("F27").Value = -8
$ws.Range("F28").Value = -9
